$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "66.552.76"
Set-TextValue "E2" "  -0.16%  "
Set-TextValue "D3" "3.236.76"
Set-TextValue "E3" "  +1.38%  "
Set-TextValue "E4" "  -0.13%  "
Set-TextValue "D5" "603.87"
Set-TextValue "E5" "  +0.25%  "
Set-TextValue "D6" "156.65"
Set-TextValue "E6" "  +0.08%  "
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "D8" "3.236.37"
Set-TextValue "E8" "  +1.38%  "
Set-TextValue "D9" "0.547"
Set-TextValue "E9" "  -0.60%  "
Set-TextValue "D10" "0.163"
Set-TextValue "E10" "  +2.68%  "
Set-TextValue "D11" "5.76"
Set-TextValue "E11" "  -2.62%  "
Set-TextValue "E12" "  -1.82%  "
Set-TextValue "D13" "0.0000272"
Set-TextValue "E13" "  +2.73%  "
Set-TextValue "D14" "38.82"
Set-TextValue "E14" "  -0.20%  "
Set-TextValue "D15" "3.761.75"
Set-TextValue "E15" "  +1.09%  "
Set-TextValue "D16" "66.585.92"
Set-TextValue "E16" "  -0.14%  "
Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.232.38"
Set-TextValue "E17" "  +1.21%  "
Set-TextValue "B18" "Polkadot"
Set-TextValue "C18" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D18" "7.28"
Set-TextValue "E18" "  -1.48%  "
Set-TextValue "E19" "  +1.37%  "
Set-TextValue "D20" "507.60"
Set-TextValue "E20" "  -1.42%  "
Set-TextValue "D21" "15.28"
Set-TextValue "E21" "  -0.63%  "
Set-TextValue "D22" "0.741"
Set-TextValue "E22" "  +0.77%  "
Set-TextValue "D23" "8.00"
Set-TextValue "E23" "  -1.60%  "
Set-TextValue "D24" "14.55"
Set-TextValue "E24" "  -2.42%  "
Set-TextValue "E25" "  +1.32%  "
Set-TextValue "D26" "0.168"
Set-TextValue "E26" "  +87.78%  "
Set-TextValue "E27" "  +0.10%  "
Set-TextValue "D28" "3.00"
Set-TextValue "E28" "  -0.46%  "
Set-TextValue "D29" "9.04"
Set-TextValue "E29" "  -2.32%  "
Set-TextValue "D30" "2.35"
Set-TextValue "E30" "  -2.45%  "
Set-TextValue "E31" "  -5.57%  "
Set-TextValue "E32" "  -1.73%  "
Set-TextValue "D33" "28.14"
Set-TextValue "E33" "  +0.17%  "
Set-TextValue "E34" "  +0.03%  "
Set-TextValue "E35" "  -4.83%  "
Set-TextValue "D36" "6.35"
Set-TextValue "E36" "  -2.86%  "
Set-TextValue "D37" "0.0₃0803"
Set-TextValue "E37" "  +17.10%  "
Set-TextValue "D38" "55.34"
Set-TextValue "E38" "  +0.85%  "
Set-TextValue "D39" "493.41"
Set-TextValue "E39" "  -3.96%  "
Set-TextValue "D40" "3.17"
Set-TextValue "E40" "  +10.08%  "
Set-TextValue "D41" "0.0421"
Set-TextValue "E41" "  -0.27%  "
Set-TextValue "E42" "  +2.51%  "
Set-TextValue "D43" "8.71"
Set-TextValue "E43" "  -1.84%  "
Set-TextValue "D44" "0.292"
Set-TextValue "E44" "  -3.93%  "
Set-TextValue "D45" "2.952.17"
Set-TextValue "E45" "  +3.10%  "
Set-TextValue "D46" "2.46"
Set-TextValue "E46" "  +0.25%  "
Set-TextValue "D47" "28.13"
Set-TextValue "E47" "  -1.30%  "
Set-TextValue "E48" "  +0.64%  "
Set-TextValue "D49" "0.119"
Set-TextValue "E49" "  +1.34%  "
Set-TextValue "E50" "  -0.04%  "
Set-TextValue "D51" "2.55"
Set-TextValue "E51" "  -3.27%  "

Write-Host "Applied cryptos update"